$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Misc. Data" column (J) header.
$ws.Cells.Item(1, 10).Value = "Misc. Data"

# The three existing data rows get a blank Misc. Data value. A plain empty
# string is treated as "clear the cell", so use a lone apostrophe (forces an
# empty, text-typed cell) and then drop the resulting quote-prefix style so
# the cell ends up as a normal empty shared-string cell.
$ws.Cells.Item(2, 10).Value = "'"
$ws.Cells.Item(2, 10).ClearFormats()
$ws.Cells.Item(3, 10).Value = "'"
$ws.Cells.Item(3, 10).ClearFormats()
$ws.Cells.Item(4, 10).Value = "'"
$ws.Cells.Item(4, 10).ClearFormats()

# The Authors column (E) values were re-saved with a field that grew slightly
# (extra padding whitespace between authors) due to the bug being fixed;
# update them to the new values.
$ws.Cells.Item(2, 5).Value = "[Luciano%Gattinoni%gattinoniluciano@gmail.com%2,                         Davide%Chiumello%NULL%3,                         Sandra%Rossi%NULL%3]"
$ws.Cells.Item(3, 5).Value = "[Luciano%Gattinoni%NULL%0,                         Silvia%Coppola%NULL%2,                         Silvia%Coppola%NULL%0,                         Massimo%Cressoni%NULL%1,                         Mattia%Busana%NULL%2,                         Mattia%Busana%NULL%0,                         Sandra%Rossi%NULL%0,                         Sandra%Rossi%NULL%0,                         Davide%Chiumello%NULL%0,                         Davide%Chiumello%NULL%0]"
$ws.Cells.Item(4, 5).Value = "[Khai%Tran%NULL%1,                         Karen%Cimon%NULL%1,                         Melissa%Severn%NULL%1,                         Carmem L.%Pessoa-Silva%NULL%1,                         John%Conly%NULL%1,                         Malcolm Gracie%Semple%NULL%2,                         Malcolm Gracie%Semple%NULL%0]"
